$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: change existing category label, add algorithm name
$ws.Range("C2").Value = "double ptr"
$ws.Range("D2").Value = "二分法"

# Row 3: add category label + new algorithm name
$ws.Range("C3").Value = "double ptr"
$ws.Range("D3").Value = "快慢指针"

# New header for the count column, plus the counts themselves
$ws.Range("E1").Value = "次数"
$ws.Range("E2").Value = 1
$ws.Range("E3").Value = 1

# Row 1 height
$ws.Rows("1").RowHeight = 24.75

# Move the active selection cursor
$ws.Range("E7").Select()
